$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy number/text formatting for the date (C) and remarks (F) columns
# from the row above (row 28) so the new row matches the existing table
# styling exactly (columns B/D/E already inherit the correct default
# column styles when a value is written into a previously-empty cell).
$ws.Range("C28").Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("F28").Copy() | Out-Null
$ws.Range("F29").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Fill in the new "Documentation" log entry (row 29)
$ws.Range("B29").Value2 = "Documentation"
$ws.Range("C29").Value2 = 43639
$ws.Range("D29").Value2 = 5
$ws.Range("E29").Value2 = 0
$ws.Range("F29").Value2 = "1) Restructuring, removing irrelevant files, creating JSDocs & documentation"

# --- Update the current selection to reflect where the author left off
$ws.Range("F30").Select()
